# Commit: "change addressbook in diagram to MeetingBook"
#
# The UndoRedo activity diagram has two shapes whose text mentions the
# old "address book" naming:
#   1. TextBox 47  -> "[command commits address book]"
#   2. Rectangle: Rounded Corners 50
#        -> "Purge redundant states and then save address book to addressBookStateList "
#
# We rename the "address book" mention that refers to the MeetingBook
# object to "MeetingBook" (shape 1), and rename the
# "addressBookStateList" identifier to "meetingBookStateList" (shape 2).
#
# We walk every shape on every slide and match on the *current* text so
# the script is resilient to shape re-ordering / re-numbering, rather
# than depending on brittle, hard-coded shape indices.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if (-not $shape.HasTextFrame) {
            continue
        }
        if (-not $shape.TextFrame.HasText) {
            continue
        }

        $tr = $shape.TextFrame.TextRange
        $text = $tr.Text

        # --- Shape 1: "...command commits address book]" -----------------
        # Split so the replacement ("MeetingBook") lands in its own run,
        # matching how PowerPoint splits runs when you select-and-retype
        # a substring.
        if ($text.Contains("command commits address book]")) {
            $needle = "address book"
            $idx = $text.IndexOf($needle)
            if ($idx -ge 0) {
                $sub = $tr.Characters($idx + 1, $needle.Length)
                $sub.Text = "MeetingBook"
            }
        }

        # --- Shape 2: "...save address book to addressBookStateList " ----
        if ($text.Contains("addressBookStateList")) {
            $needle2 = "addressBookStateList"
            $idx2 = $text.IndexOf($needle2)
            if ($idx2 -ge 0) {
                $sub2 = $tr.Characters($idx2 + 1, $needle2.Length)
                $sub2.Text = "meetingBookStateList"
            }
        }
    }
}
